# Fruta / hortaliza, semanal
# Insert two new weekly price-report rows (2021-11-16) above the existing
# "Haba" records in the consolidated sheet, pushing the former rows
# 176-186 down to 178-188.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 176 and 177 (existing data shifts down).
$ws.Rows.Item(176).Insert()
$ws.Rows.Item(177).Insert()

# --- New row 176 ---
$ws.Cells.Item(176, 1).Value = 6
$ws.Cells.Item(176, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(176, 3).Value = "Metropolitana"
$ws.Cells.Item(176, 4).Value = 44516
$ws.Cells.Item(176, 5).Value = 13
$ws.Cells.Item(176, 6).Value = 100112026
$ws.Cells.Item(176, 7).Value = "Haba"
$ws.Cells.Item(176, 8).Value = "Sin especificar"
$ws.Cells.Item(176, 9).Value = "Primera"
$ws.Cells.Item(176, 10).Value = 350
$ws.Cells.Item(176, 11).Value = 5000
$ws.Cells.Item(176, 12).Value = 6000
$ws.Cells.Item(176, 13).Value = 5657
$ws.Cells.Item(176, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(176, 15).Value = "Región Metropolitana"
$ws.Cells.Item(176, 16).Value = 226
$ws.Cells.Item(176, 17).Value = 25
$ws.Cells.Item(176, 18).Value = "Hortaliza"

# --- New row 177 ---
$ws.Cells.Item(177, 1).Value = 6
$ws.Cells.Item(177, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(177, 3).Value = "Metropolitana"
$ws.Cells.Item(177, 4).Value = 44516
$ws.Cells.Item(177, 5).Value = 13
$ws.Cells.Item(177, 6).Value = 100112026
$ws.Cells.Item(177, 7).Value = "Haba"
$ws.Cells.Item(177, 8).Value = "Sin especificar"
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 1010
$ws.Cells.Item(177, 11).Value = 5000
$ws.Cells.Item(177, 12).Value = 7000
$ws.Cells.Item(177, 13).Value = 6178
$ws.Cells.Item(177, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(177, 15).Value = "Región del Maule"
$ws.Cells.Item(177, 16).Value = 247
$ws.Cells.Item(177, 17).Value = 25
$ws.Cells.Item(177, 18).Value = "Hortaliza"
